$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("ForecastValidations")
$wsCurrent  = $wb.Worksheets.Item("CurrentWeatherValidations")

# --- ForecastValidations content edits -----------------------------------

# B3: "Aundh Camp" -> "AundhCamp"
$wsForecast.Range("B3").Value = "AundhCamp"

# B4: "New Delhi" -> "NewDelhi"
$wsForecast.Range("B4").Value = "NewDelhi"

# D4: text " 77.2311" -> numeric 77.2311
$wsForecast.Range("D4").Value = 77.2311

# B2 picks up the highlight style already used on CurrentWeatherValidations!B2
$wsCurrent.Range("B2").Copy()
$wsForecast.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Selections / active sheet -------------------------------------------

[void]$wsCurrent.Range("B2").Select()
$wsCurrent.Activate()

[void]$wsForecast.Range("D16").Select()
$wsForecast.Activate()
